# Mise à jour de l'application
# Adds the week of 2026-01-27 (serial date 46049) training-load entries
# (rows 762-773) to the "Feuil1" sheet, mirroring the style/formatting of
# the rows immediately above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 762
$lastNewRow  = 773
$nbsp = [char]0x00A0

# 1) Stamp the formatting of the last existing data row (761) onto the
#    12 new rows so number formats / fonts match the rest of the table.
$ws.Range("A761:I761").Copy()
$ws.Range("A" + $firstNewRow + ":I" + $lastNewRow).PasteSpecial(-4122)

# 2) Row data: date, player, volume, intensite, charge, douleur,
#    localisation douleur, plaisir
$rows = @(
    @{ Row=762; Date=46049; Joueur="Amir Etien";      Volume=5; Intensite=3; Charge=6; Douleur=0; Loc="Genou";      Plaisir=0 },
    @{ Row=763; Date=46049; Joueur="Yoan Zouma";       Volume=5; Intensite=5; Charge=0; Douleur=0; Loc="";          Plaisir=0 },
    @{ Row=764; Date=46049; Joueur="Theo Owono";       Volume=3; Intensite=5; Charge=0; Douleur=10; Loc="";         Plaisir=10 },
    @{ Row=765; Date=46049; Joueur="Malik Boussaid";   Volume=6; Intensite=7; Charge=0; Douleur=8; Loc="";          Plaisir=8 },
    @{ Row=766; Date=46049; Joueur="Jeremie Laurent";  Volume=6; Intensite=6; Charge=0; Douleur=9; Loc="";          Plaisir=9 },
    @{ Row=767; Date=46049; Joueur="Ilan Ihaddadene";  Volume=6; Intensite=0; Charge=9; Douleur=9; Loc="";          Plaisir=9 },
    @{ Row=768; Date=46049; Joueur="Romain Thunet";    Volume=7; Intensite=4; Charge=2; Douleur=7; Loc="Terrain";   Plaisir=7 },
    @{ Row=769; Date=46049; Joueur="Naim Dhib";        Volume=5; Intensite=3; Charge=3; Douleur=3; Loc=("Psoas" + $nbsp); Plaisir=3 },
    @{ Row=770; Date=46049; Joueur="Hedi Nasri";       Volume=5; Intensite=5; Charge=3; Douleur=6; Loc=("Hanche" + $nbsp); Plaisir=6 },
    @{ Row=771; Date=46049; Joueur="Naim Ighbane";     Volume=6; Intensite=6; Charge=3; Douleur=9; Loc="Coup tibia"; Plaisir=9 },
    @{ Row=772; Date=46049; Joueur="Karahali Souaré";  Volume=6; Intensite=7; Charge=6; Douleur=10; Loc="Cheville"; Plaisir=10 },
    @{ Row=773; Date=46049; Joueur="Mattheo Haon";     Volume=7; Intensite=8; Charge=0; Douleur=4; Loc="";          Plaisir=4 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.Date
    $ws.Range("B$n").Value = $r.Joueur
    $ws.Range("C$n").Value = 70
    $ws.Range("D$n").Value = $r.Volume
    $ws.Range("E$n").Value = $r.Intensite
    $ws.Range("F$n").Value = $r.Charge
    if ($r.Loc -ne "") {
        $ws.Range("G$n").Value = $r.Loc
    } else {
        # keep the same "empty but formatted" look as the rest of the sheet
        $ws.Range("G760").Copy()
        $ws.Range("G$n").PasteSpecial(-4122)
    }
    $ws.Range("H$n").Value = $r.Douleur
}

# 3) Formulas for column I (Charge = Volume * Intensite), filled in two
#    batches so the split mirrors how the data was actually entered.
$ws.Range("I762:I771").Formula = "=C762*D762"
$ws.Range("I772:I773").Formula = "=C772*D772"

# 4) Leave the sheet scrolled/selected near the newly entered rows, like
#    the author would after typing the last row.
$ws.Range("C776").Select()

Write-Output "Added rows $firstNewRow to $lastNewRow (week of 2026-01-27)."
